$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - reorder "Recorded By" list
$ws.Range("G2").Value = "System, servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

# Row 3 - reorder "Recorded By" list
$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# Row 4 - reorder "Recorded By" list
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

# Row 5 - reorder "Recorded By" list
$ws.Range("G5").Value = "Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# Row 7 - new recorder added and reordered list, attendance count updated
$ws.Range("G7").Value = "lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Amera.a.saad@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"
$ws.Range("H7").Value = "79/251"

# Row 10 - updated average attendance percentage for HISTOLOGY metric
# (write as literal text, then restore the original cell formatting so
# Excel's auto percent-to-number conversion doesn't change the cell's
# style/number format)
$ws.Range("L10").Value = "'24.3%"
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 12 - reorder "Recorded By" list
$ws.Range("G12").Value = "dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg"

# Row 15 - updated group average attendance percentage
$ws.Range("S15").Value = "'24.3%"
$ws.Range("L10").Copy()
$ws.Range("S15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
